$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Npy / Prlhr / MuSCs -> "ECs") keeps the same shared-string index for
# its Target cluster cell (D2), but the updated TPM pipeline re-ran the
# dedup of the shared-strings table so that slot now resolves to "FAPs".
$ws.Range("D2").Value = "FAPs"

# Row 2 gets the expression/specificity numbers that used to belong to the
# (now removed) row 3.
$ws.Range("M2").Value = 0.001793666666666667
$ws.Range("N2").Value = 0.005381
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0001068248077777778
$ws.Range("R2").Value = 0.00096142327
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Row 3 (the old MuSCs -> FAPs pair) is dropped entirely.
$ws.Rows.Item(3).Delete()
